# Applies the "Added Graphs and Ray trace for Real lenses" edit:
#  - Update view state (scroll position / active selection) on Sheet1
#  - Swap/update numeric values for M3, N3 (EFL1 params) and AO3 (PETC for lens 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the cell values
$ws.Range("M3").Value = 1000000000
$ws.Range("N3").Value = -17.5
$ws.Range("AO3").Value = 8.9700000000000006

# Update the view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 31   # Column AE = 31
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AL6").Select()
